$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 2 (primary header) - BTec_Logo-Orange picture: image1.jpg -> image2.jpg
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists) {
    $r = $hdr.Range
    for ($i = 1; $i -le $r.InlineShapes.Count; $i++) {
        $sh = $r.InlineShapes.Item($i)
        if ($sh.AlternativeText -eq "BTec_Logo-Orange") {
            $sh.Name = "image2.jpg"
        }
    }
}

# Footer 1 - Pearson logo picture: image2.png -> image1.png
$ftr1 = $sec.Footers.Item(1)
if ($ftr1.Exists) {
    $r = $ftr1.Range
    for ($i = 1; $i -le $r.InlineShapes.Count; $i++) {
        $sh = $r.InlineShapes.Item($i)
        if ($sh.AlternativeText -like "*PearsonLogo.png") {
            $sh.Name = "image1.png"
        }
    }
}

# Footer 2 - Pearson logo picture: image2.png -> image1.png
$ftr2 = $sec.Footers.Item(2)
if ($ftr2.Exists) {
    $r = $ftr2.Range
    for ($i = 1; $i -le $r.InlineShapes.Count; $i++) {
        $sh = $r.InlineShapes.Item($i)
        if ($sh.AlternativeText -like "*PearsonLogo.png") {
            $sh.Name = "image1.png"
        }
    }
}
